# Add "Planning court" sheet to the London Administrative Court daily cause
# list: duplicate the existing sheet so the two courts have identical
# columns/data, then rename the tabs and restore the correct selections.

$wb = $excel.ActiveWorkbook
$original = $wb.Worksheets.Item(1)

# Copy the original sheet, placing the copy immediately before it. This
# becomes the new, active "Hearing list" tab; the original sheet (now
# second) becomes "Planning court".
$original.Copy($original)

$hearingList = $wb.Worksheets.Item(1)
$planningCourt = $wb.Worksheets.Item(2)

$hearingList.Name = "Hearing list"
$planningCourt.Name = "Planning court"

# "Hearing list" is the active/selected tab, with F33 as the selected cell.
$hearingList.Activate()
[void]$hearingList.Range("F33").Select()
